# [MOSIP-15491] Machine onboarding code made compatible with 1.1.5.4.
#
# Changes applied:
#  1. Shared string "Reg-Mach-1" -> "desktop-fvgt677" (used by cells A2 & A3
#     on Sheet1, the sample machine-name rows).
#  2. Active-cell selection moved from A2 to C3.
#  3. Minor column-width retouches on columns E (5) and H (8) to match the
#     re-saved sheet metrics (closest values reachable through the
#     ColumnWidth object model, which is quantised to 1/6 of a character).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the machine name text wherever it appears (A2:A3 share the
#    "Reg-Mach-1" string in the shared-string table).
$ws.Range("A2:A3").Value = "desktop-fvgt677"

# 3. Column width tweaks (E -> 38.62, H -> 10.11 characters in the saved
#    OOXML "width" attribute, which is ColumnWidth + 5/6).
$ws.Columns.Item(5).ColumnWidth = 37.78666666666667
$ws.Columns.Item(8).ColumnWidth = 9.276666666666667

# 2. Move the active cell / selection to C3.
$ws.Range("C3").Select()
